$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 10.36733066666667
$ws.Range("H2").Value2 = 31.101992
$ws.Range("I2").Value2 = 0.1169328841728879
$ws.Range("J2").Value2 = 0.1169328841728879
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 7.248785666666667
$ws.Range("N2").Value2 = 21.746357
$ws.Range("O2").Value2 = 0.07891374419744837
$ws.Range("P2").Value2 = 0.07891374419744837
$ws.Range("Q2").Value2 = 75.15055793812711
$ws.Range("R2").Value2 = 676.355021443144
$ws.Range("S2").Value2 = 0.009227611709889135
$ws.Range("T2").Value2 = 0.009227611709889137

$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 10.36733066666667
$ws.Range("H3").Value2 = 31.101992
$ws.Range("I3").Value2 = 0.1169328841728879
$ws.Range("J3").Value2 = 0.1169328841728879
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 37.82684066666667
$ws.Range("N3").Value2 = 113.480522
$ws.Range("O3").Value2 = 0.4118010609547572
$ws.Range("P3").Value2 = 0.4118010609547572
$ws.Range("Q3").Value2 = 392.1633652666471
$ws.Range("R3").Value2 = 3529.470287399824
$ws.Range("S3").Value2 = 0.04815308576289498
$ws.Range("T3").Value2 = 0.04815308576289499

$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 10.36733066666667
$ws.Range("H4").Value2 = 31.101992
$ws.Range("I4").Value2 = 0.1169328841728879
$ws.Range("J4").Value2 = 0.1169328841728879
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 38.20927633333334
$ws.Range("N4").Value2 = 114.627829
$ws.Range("O4").Value2 = 0.4159644383477588
$ws.Range("P4").Value2 = 0.4159644383477588
$ws.Range("Q4").Value2 = 396.1282022817076
$ws.Range("R4").Value2 = 3565.153820535368
$ws.Range("S4").Value2 = 0.04863992148935885
$ws.Range("T4").Value2 = 0.04863992148935886

$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 10.36733066666667
$ws.Range("H5").Value2 = 31.101992
$ws.Range("I5").Value2 = 0.1169328841728879
$ws.Range("J5").Value2 = 0.1169328841728879
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 8.572171666666666
$ws.Range("N5").Value2 = 25.716515
$ws.Range("O5").Value2 = 0.09332075650003555
$ws.Range("P5").Value2 = 0.09332075650003555
$ws.Range("Q5").Value2 = 88.87053819976443
$ws.Range("R5").Value2 = 799.8348437978798
$ws.Range("S5").Value2 = 0.01091226521074493
$ws.Range("T5").Value2 = 0.01091226521074493

$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 37.91490933333333
$ws.Range("H6").Value2 = 113.744728
$ws.Range("I6").Value2 = 0.4276413904453658
$ws.Range("J6").Value2 = 0.4276413904453659
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 7.248785666666667
$ws.Range("N6").Value2 = 21.746357
$ws.Range("O6").Value2 = 0.07891374419744837
$ws.Range("P6").Value2 = 0.07891374419744837
$ws.Range("Q6").Value2 = 274.8370513284329
$ws.Range("R6").Value2 = 2473.533461955897
$ws.Range("S6").Value2 = 0.03374678329384673
$ws.Range("T6").Value2 = 0.03374678329384674

$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 37.91490933333333
$ws.Range("H7").Value2 = 113.744728
$ws.Range("I7").Value2 = 0.4276413904453658
$ws.Range("J7").Value2 = 0.4276413904453659
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 37.82684066666667
$ws.Range("N7").Value2 = 113.480522
$ws.Range("O7").Value2 = 0.4118010609547572
$ws.Range("P7").Value2 = 0.4118010609547572
$ws.Range("Q7").Value2 = 1434.201234243113
$ws.Range("R7").Value2 = 12907.81110818802
$ws.Range("S7").Value2 = 0.1761031782935692
$ws.Range("T7").Value2 = 0.1761031782935693

$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 37.91490933333333
$ws.Range("H8").Value2 = 113.744728
$ws.Range("I8").Value2 = 0.4276413904453658
$ws.Range("J8").Value2 = 0.4276413904453659
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 38.20927633333334
$ws.Range("N8").Value2 = 114.627829
$ws.Range("O8").Value2 = 0.4159644383477588
$ws.Range("P8").Value2 = 0.4159644383477588
$ws.Range("Q8").Value2 = 1448.701247870613
$ws.Range("R8").Value2 = 13038.31123083551
$ws.Range("S8").Value2 = 0.1778836107908612
$ws.Range("T8").Value2 = 0.1778836107908612

$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 37.91490933333333
$ws.Range("H9").Value2 = 113.744728
$ws.Range("I9").Value2 = 0.4276413904453658
$ws.Range("J9").Value2 = 0.4276413904453659
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 8.572171666666666
$ws.Range("N9").Value2 = 25.716515
$ws.Range("O9").Value2 = 0.09332075650003555
$ws.Range("P9").Value2 = 0.09332075650003555
$ws.Range("Q9").Value2 = 325.0131115314355
$ws.Range("R9").Value2 = 2925.11800378292
$ws.Range("S9").Value2 = 0.0399078180670886
$ws.Range("T9").Value2 = 0.03990781806708862

$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 26.72147866666667
$ws.Range("H10").Value2 = 80.164436
$ws.Range("I10").Value2 = 0.3013909433702152
$ws.Range("J10").Value2 = 0.3013909433702153
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 7.248785666666667
$ws.Range("N10").Value2 = 21.746357
$ws.Range("O10").Value2 = 0.07891374419744837
$ws.Range("P10").Value2 = 0.07891374419744837
$ws.Range("Q10").Value2 = 193.6982715510725
$ws.Range("R10").Value2 = 1743.284443959652
$ws.Range("S10").Value2 = 0.02378388780854481
$ws.Range("T10").Value2 = 0.02378388780854481

$ws.Range("E11").Value2 = 3
$ws.Range("G11").Value2 = 26.72147866666667
$ws.Range("H11").Value2 = 80.164436
$ws.Range("I11").Value2 = 0.3013909433702152
$ws.Range("J11").Value2 = 0.3013909433702153
$ws.Range("K11").Value2 = 3
$ws.Range("M11").Value2 = 37.82684066666667
$ws.Range("N11").Value2 = 113.480522
$ws.Range("O11").Value2 = 0.4118010609547572
$ws.Range("P11").Value2 = 0.4118010609547572
$ws.Range("Q11").Value2 = 1010.789115901733
$ws.Range("R11").Value2 = 9097.102043115592
$ws.Range("S11").Value2 = 0.1241131102420098
$ws.Range("T11").Value2 = 0.1241131102420098

$ws.Range("E12").Value2 = 3
$ws.Range("G12").Value2 = 26.72147866666667
$ws.Range("H12").Value2 = 80.164436
$ws.Range("I12").Value2 = 0.3013909433702152
$ws.Range("J12").Value2 = 0.3013909433702153
$ws.Range("K12").Value2 = 3
$ws.Range("M12").Value2 = 38.20927633333334
$ws.Range("N12").Value2 = 114.627829
$ws.Range("O12").Value2 = 0.4159644383477588
$ws.Range("P12").Value2 = 0.4159644383477588
$ws.Range("Q12").Value2 = 1021.008362409938
$ws.Range("R12").Value2 = 9189.075261689444
$ws.Range("S12").Value2 = 0.1253679144820928
$ws.Range("T12").Value2 = 0.1253679144820928

$ws.Range("E13").Value2 = 3
$ws.Range("G13").Value2 = 26.72147866666667
$ws.Range("H13").Value2 = 80.164436
$ws.Range("I13").Value2 = 0.3013909433702152
$ws.Range("J13").Value2 = 0.3013909433702153
$ws.Range("K13").Value2 = 3
$ws.Range("M13").Value2 = 8.572171666666666
$ws.Range("N13").Value2 = 25.716515
$ws.Range("O13").Value2 = 0.09332075650003555
$ws.Range("P13").Value2 = 0.09332075650003555
$ws.Range("Q13").Value2 = 229.0611023178378
$ws.Range("R13").Value2 = 2061.54992086054
$ws.Range("S13").Value2 = 0.02812603083756786
$ws.Range("T13").Value2 = 0.02812603083756786

$ws.Range("E14").Value2 = 3
$ws.Range("G14").Value2 = 13.65680433333333
$ws.Range("H14").Value2 = 40.970413
$ws.Range("I14").Value2 = 0.154034782011531
$ws.Range("J14").Value2 = 0.154034782011531
$ws.Range("K14").Value2 = 3
$ws.Range("M14").Value2 = 7.248785666666667
$ws.Range("N14").Value2 = 21.746357
$ws.Range("O14").Value2 = 0.07891374419744837
$ws.Range("P14").Value2 = 0.07891374419744837
$ws.Range("Q14").Value2 = 98.9952475039379
$ws.Range("R14").Value2 = 890.9572275354411
$ws.Range("S14").Value2 = 0.01215546138516768
$ws.Range("T14").Value2 = 0.01215546138516768

$ws.Range("E15").Value2 = 3
$ws.Range("G15").Value2 = 13.65680433333333
$ws.Range("H15").Value2 = 40.970413
$ws.Range("I15").Value2 = 0.154034782011531
$ws.Range("J15").Value2 = 0.154034782011531
$ws.Range("K15").Value2 = 3
$ws.Range("M15").Value2 = 37.82684066666667
$ws.Range("N15").Value2 = 113.480522
$ws.Range("O15").Value2 = 0.4118010609547572
$ws.Range("P15").Value2 = 0.4118010609547572
$ws.Range("Q15").Value2 = 516.5937615328429
$ws.Range("R15").Value2 = 4649.343853795586
$ws.Range("S15").Value2 = 0.06343168665628322
$ws.Range("T15").Value2 = 0.06343168665628324

$ws.Range("E16").Value2 = 3
$ws.Range("G16").Value2 = 13.65680433333333
$ws.Range("H16").Value2 = 40.970413
$ws.Range("I16").Value2 = 0.154034782011531
$ws.Range("J16").Value2 = 0.154034782011531
$ws.Range("K16").Value2 = 3
$ws.Range("M16").Value2 = 38.20927633333334
$ws.Range("N16").Value2 = 114.627829
$ws.Range("O16").Value2 = 0.4159644383477588
$ws.Range("P16").Value2 = 0.4159644383477588
$ws.Range("Q16").Value2 = 521.8166106025975
$ws.Range("R16").Value2 = 4696.349495423377
$ws.Range("S16").Value2 = 0.06407299158544595
$ws.Range("T16").Value2 = 0.06407299158544597

$ws.Range("E17").Value2 = 3
$ws.Range("G17").Value2 = 13.65680433333333
$ws.Range("H17").Value2 = 40.970413
$ws.Range("I17").Value2 = 0.154034782011531
$ws.Range("J17").Value2 = 0.154034782011531
$ws.Range("K17").Value2 = 3
$ws.Range("M17").Value2 = 8.572171666666666
$ws.Range("N17").Value2 = 25.716515
$ws.Range("O17").Value2 = 0.09332075650003555
$ws.Range("P17").Value2 = 0.09332075650003555
$ws.Range("Q17").Value2 = 117.0684711634105
$ws.Range("R17").Value2 = 1053.616240470695
$ws.Range("S17").Value2 = 0.01437464238463414
$ws.Range("T17").Value2 = 0.01437464238463414
